$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Wed May 31 22:21:44 UTC 2023 with GitHub Actions

$ws.Range("D2").Value = '27.129.45'
$ws.Range("E2").Value = '  -2.07%  '
$ws.Range("D3").Value = '1.867.41'
$ws.Range("E3").Value = '  -1.93%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.09'
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5148'
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3759'
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07145'
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8908'
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.74'
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.894.31'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07550'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.309'
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.59'
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008487'
$ws.Range("E17").Value = '  -2.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.10'
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '27.168.40'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.994'
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("D22").Value = '2.091.49'
$ws.Range("E22").Value = '  -2.08%  '
$ws.Range("E23").Value = '  -3.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.453'
$ws.Range("E24").Value = '  -2.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.835'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.01'
$ws.Range("E26").Value = '  -5.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.94'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.089'
$ws.Range("E28").Value = '  -3.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.71'
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.660'
$ws.Range("E30").Value = '  -3.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.679'
$ws.Range("E31").Value = '  -3.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09225'
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05112'
$ws.Range("E33").Value = '  -3.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.076'
$ws.Range("E34").Value = '  -3.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.160'
$ws.Range("E35").Value = '  -5.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7248'
$ws.Range("E36").Value = '  -7.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02037'
$ws.Range("E37").Value = '  -2.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.097'
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.503'
$ws.Range("E39").Value = '  -3.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.077'
$ws.Range("E40").Value = '  -1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5290'
$ws.Range("E41").Value = '  -4.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.504'
$ws.Range("E42").Value = '  -3.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.64'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.330'
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1470'
$ws.Range("E45").Value = '  -3.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9998'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4620'
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.957'
$ws.Range("E48").Value = '  -5.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.559'
$ws.Range("E49").Value = '  -3.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.64'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.59'
$ws.Range("E51").Value = '  -4.66%  '
